# Update ticker price timestamps and values (commit: "Add bat and update time")

$wb = $excel.ActiveWorkbook

# --- btc_ticker_price ---
$ws = $wb.Worksheets.Item("btc_ticker_price")
$ws.Range("B2").Value = 1526350746056
$ws.Range("B3").Value = 8673.370000000001
$ws.Range("C3").Value = 8669.6636
$ws.Range("D3").Value = 8664.42
$ws.Range("B4").Value = 8679.860000000001
$ws.Range("C4").Value = 8687.1564
$ws.Range("D4").Value = 8676.48

# --- eth_ticker_price ---
$ws = $wb.Worksheets.Item("eth_ticker_price")
$ws.Range("B2").Value = 1526350746080
$ws.Range("B3").Value = 726.21
$ws.Range("C3").Value = 725.86
$ws.Range("D3").Value = 726.6799999999999
$ws.Range("B4").Value = 726.27
$ws.Range("C4").Value = 726.942
$ws.Range("D4").Value = 727.34

# --- usdt_c2c_price ---
$ws = $wb.Worksheets.Item("usdt_c2c_price")
$ws.Range("B2").Value = 1526350750988
$ws.Range("B3").Value = 6.5
$ws.Range("B4").Value = 6.51

# --- vhkd_c2c_price ---
$ws = $wb.Worksheets.Item("vhkd_c2c_price")
$ws.Range("B2").Value = 1526350751119
$ws.Range("B4").Value = 0.8288

# --- btc_c2c_price ---
$ws = $wb.Worksheets.Item("btc_c2c_price")
$ws.Range("B2").Value = 1526350751704
$ws.Range("B3").Value = 56412
$ws.Range("B4").Value = 56500

# --- eth_c2c_price ---
$ws = $wb.Worksheets.Item("eth_c2c_price")
$ws.Range("B2").Value = 1526350752190
$ws.Range("B3").Value = 4733
$ws.Range("B4").Value = 4740
